$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Delete the "Meta description: ..." paragraph that currently
#        follows the title heading near the top of the document. ---
$metaPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Meta description:*") {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    [void]$metaPara.Range.Delete()
}

# --- 2. Turn the final paragraph (previously the italic image-generation
#        prompt) into two paragraphs: a new bold paragraph carrying the
#        page title, followed by the (still italic) meta-description
#        sentence that used to live in the paragraph removed above. ---
$lastPara = $d.Paragraphs.Last
$rng = $lastPara.Range
$rng.Collapse(1)
$xml = "<w:p $wns><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Drift King Free: Review of High-Octane Racing Slot Game</w:t></w:r></w:p><w:p $wns><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Take a look at our review of Drift King, a high-octane racing slot game with exciting bonus features. Play now for free!</w:t></w:r></w:p>"
[void]$rng.InsertXML($xml)
